$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells F1:I1
$ws.Range("F1").Value = "id"
$ws.Range("G1").Value = "source_file"
$ws.Range("H1").Value = "text"
$ws.Range("I1").Value = "review_result"

# Copy header style from an existing header cell (A1) to the new headers
$ws.Range("A1").Copy()
$ws.Range("F1:I1").PasteSpecial(-4122)  # xlPasteFormats

# New data row 2
$ws.Range("A2").Value = "parisk"
$ws.Range("B2").Value = 3
$ws.Range("C2").Value = "nan"
$ws.Range("D2").Value = "DIS"
$ws.Range("E2").Value = "WRI"
$ws.Range("F2").Value = "1269f1fb-9c21-42a9-ae5e-c80f92622adc"
$ws.Range("G2").Value = "Bk6qQGWRb_annotated.xlsx"
$ws.Range("H2").Value = "Then how bootstrap dqn extend the idea to deep learning, followed by the noisy net, bbq, shallow UBE and LS-DQN."
$ws.Range("I2").Value = "Correct"
